$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "Dr. Gehan Adel, Administrator, Dr. Servinaz Sayed Mohammad, Dr. Veronia Rafat, Dr. Amira Sobhy"
$ws.Range("G3").Value = "Dr. Eman Tantawi, Administrator, Dr. Asmaa Reda, Dr. Majorelle Magdy, Dr. Hend Mahmoud, Dr. Veronia Rafat"
$ws.Range("G4").Value = "Dr. Eman Tantawi, Dr. Gehan Adel, Dr. Asmaa Reda, Dr. Hend Mahmoud, Dr. Majorelle Magdy, Dr. Servinaz Sayed Mohammad"
$ws.Range("G5").Value = "Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Eman Tantawi, Dr. Asmaa Reda"
$ws.Range("G6").Value = "Dr. Manar Montaser, Dr. Majorelle Magdy, Dr. Alshimaa Atef, Dr. Mohammad El-Tanany, Dr. Menna tuâ€™Allah Medhat"
$ws.Range("G7").Value = "Dr. Nada Mohammad, Dr. Menna tu'Alllah Mohammad, Dr. Abeer Ragab, Dr. Lamiaa Ossama, Dr. Fatma Elhady, Dr. Amera Ahmad Saad, Dr. Kerelos Zareef"
$ws.Range("G8").Value = "Dr. Abeer Ragab, Dr. Nada Mohammad"
$ws.Range("G11").Value = "Dr. Aya Saeed, Dr. Amal Awwad, Dr. Safa Hany"
$ws.Range("G12").Value = "Dr. Marina Youhanna, Dr. Yasmeena Fattoh, Dr. Eman M. Abo-Sakaya, Dr. Dina Adel, Dr. Madeha Saeed, Dr. Amira Ibrahim"
$ws.Range("G13").Value = "Dr. Amira Ibrahim, Dr. Yasmeena Fattoh, Dr. Esraa Mostafa"
$ws.Range("G15").Value = "Dr. Mohammad Safwat, Dr. Rania Ahmad Youssef"
$ws.Range("G17").Value = "Dr. Mohammad Safwat, Dr. Esraa Samy"
$ws.Range("G28").Value = "Dr. Maryam Ashraf, Dr. Aya Emad"
$ws.Range("G30").Value = "Dr. Shorok Mohammad, Dr. Aya Hanafy, Dr. Yassmen Ahmad, Dr. Wafaa Ebida"
